# Update market-price derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# in the per-job Leve profit tables, refreshed from the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 3328.111
$ws.Range("I62").Value = 3007.8572
$ws.Range("J62").Value = 4449
$ws.Range("K62").Value = 3007.8572
$ws.Range("L62").Value = 4449
$ws.Range("M62").Value = -2383.8572
$ws.Range("N62").Value = -5697

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 3328.111
$ws.Range("I65").Value = 3007.8572
$ws.Range("J65").Value = 4449
$ws.Range("K65").Value = 15039.286
$ws.Range("L65").Value = 22245
$ws.Range("M65").Value = -11919.286
$ws.Range("N65").Value = -28485

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 4167.125
$ws.Range("I86").Value = 2379.4546
$ws.Range("K86").Value = 2379.4546
$ws.Range("M86").Value = -1256.4546

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 4167.125
$ws.Range("I89").Value = 2379.4546
$ws.Range("K89").Value = 11897.273
$ws.Range("M89").Value = -6281.273000000001

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 100001460
$ws.Range("I113").Value = 100001460
$ws.Range("K113").Value = 100001460
$ws.Range("M113").Value = -99998206

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 6527.5
$ws.Range("I2").Value = 907
$ws.Range("J2").Value = 11023.9
$ws.Range("K2").Value = 907
$ws.Range("L2").Value = 11023.9
$ws.Range("M2").Value = -794
$ws.Range("N2").Value = -11249.9

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2510.2222
$ws.Range("I32").Value = 2641.0327
$ws.Range("K32").Value = 2641.0327
$ws.Range("M32").Value = -2354.0327

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 6527.5
$ws.Range("I116").Value = 907
$ws.Range("J116").Value = 11023.9
$ws.Range("K116").Value = 907
$ws.Range("L116").Value = 11023.9
$ws.Range("M116").Value = 1387
$ws.Range("N116").Value = -15611.9

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1108.7778
$ws.Range("I122").Value = 1120.8235
$ws.Range("K122").Value = 3362.4705
$ws.Range("M122").Value = -912.4704999999999

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2335.568
$ws.Range("I132").Value = 2072.4473
$ws.Range("K132").Value = 6217.341899999999
$ws.Range("M132").Value = -3687.341899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2526.5217
$ws.Range("I86").Value = 2827.0344
$ws.Range("J86").Value = 2013.8823
$ws.Range("K86").Value = 2827.0344
$ws.Range("L86").Value = 2013.8823
$ws.Range("M86").Value = -1704.0344
$ws.Range("N86").Value = -4259.8823

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2526.5217
$ws.Range("I89").Value = 2827.0344
$ws.Range("J89").Value = 2013.8823
$ws.Range("K89").Value = 14135.172
$ws.Range("L89").Value = 10069.4115
$ws.Range("M89").Value = -8519.172
$ws.Range("N89").Value = -21301.4115

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1468.5454
$ws.Range("I31").Value = 910.7143
$ws.Range("K31").Value = 910.7143
$ws.Range("M31").Value = -615.7143

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1468.5454
$ws.Range("I34").Value = 910.7143
$ws.Range("K34").Value = 910.7143
$ws.Range("M34").Value = -708.7143

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 12469.7
$ws.Range("I132").Value = 16554.715
$ws.Range("K132").Value = 49664.145
$ws.Range("M132").Value = -47134.145

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 2329.027
$ws.Range("J68").Value = 2379.8333
$ws.Range("L68").Value = 7139.499899999999
$ws.Range("N68").Value = -8761.499899999999

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 2329.027
$ws.Range("J71").Value = 2379.8333
$ws.Range("L71").Value = 21418.4997
$ws.Range("N71").Value = -29530.4997

# Row 96: Hunger Is No Game / Popoto Soba
$ws.Range("H96").Value = 7905.08
$ws.Range("J96").Value = 7905.08
$ws.Range("L96").Value = 23715.24
$ws.Range("N96").Value = -27833.24

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me / Bone Hora
$ws.Range("H5").Value = 13577.692
$ws.Range("J5").Value = 13577.692
$ws.Range("L5").Value = 13577.692
$ws.Range("N5").Value = -13801.692

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 4626.6
$ws.Range("I80").Value = 3374.875
$ws.Range("J80").Value = 6057.143
$ws.Range("K80").Value = 3374.875
$ws.Range("L80").Value = 6057.143
$ws.Range("M80").Value = -2376.875
$ws.Range("N80").Value = -8053.143

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 4626.6
$ws.Range("I83").Value = 3374.875
$ws.Range("J83").Value = 6057.143
$ws.Range("K83").Value = 16874.375
$ws.Range("L83").Value = 30285.715
$ws.Range("M83").Value = -11882.375
$ws.Range("N83").Value = -40269.715

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 3449.111
$ws.Range("I102").Value = 6520.6665
$ws.Range("J102").Value = 1913.3334
$ws.Range("K102").Value = 6520.6665
$ws.Range("L102").Value = 1913.3334
$ws.Range("M102").Value = -4898.6665
$ws.Range("N102").Value = -5157.3334

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3453.4
$ws.Range("I122").Value = 2422.3333
$ws.Range("K122").Value = 7266.999899999999
$ws.Range("M122").Value = -4816.999899999999

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3274.0527
$ws.Range("I126").Value = 2021.1428
$ws.Range("J126").Value = 4004.9167
$ws.Range("K126").Value = 6063.428400000001
$ws.Range("L126").Value = 12014.7501
$ws.Range("M126").Value = -3593.428400000001
$ws.Range("N126").Value = -16954.7501

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2357.6
$ws.Range("I40").Value = 2357.6
$ws.Range("K40").Value = 2357.6
$ws.Range("M40").Value = -2221.6

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 2166.3333
$ws.Range("I82").Value = 1499
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 1499
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1138
$ws.Range("N82").Value = -3222

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 2166.3333
$ws.Range("I85").Value = 1499
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 1499
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -251
$ws.Range("N85").Value = -4996

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 20245990
$ws.Range("I122").Value = 20245990
$ws.Range("K122").Value = 60737970
$ws.Range("M122").Value = -60735520

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables / Hempen Underpants
$ws.Range("H2").Value = 3500
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4888

# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 35722620
$ws.Range("I62").Value = 62505976
$ws.Range("J62").Value = 11483
$ws.Range("K62").Value = 62505976
$ws.Range("L62").Value = 11483
$ws.Range("M62").Value = -62505352
$ws.Range("N62").Value = -12731

# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 35722620
$ws.Range("I65").Value = 62505976
$ws.Range("J65").Value = 11483
$ws.Range("K65").Value = 312529880
$ws.Range("L65").Value = 57415
$ws.Range("M65").Value = -63655
